$d = $word.ActiveDocument

$replacements = @(
    @("87×18=1566", "38×78=2964"),
    @("70×60=4200", "68×42=2856"),
    @("78×73=5694", "13×46=598"),
    @("27×51=1377", "34×85=2890"),
    @("38×59=2242", "89×98=8722"),
    @("57×77=4389", "28×68=1904"),
    @("71×16=1136", "66×12=792"),
    @("59×60=3540", "58×51=2958"),
    @("95×53=5035", "55×49=2695"),
    @("45×81=3645", "87×99=8613"),
    @("31×43=1333", "75×21=1575"),
    @("39×78=3042", "98×50=4900"),
    @("30×55=1650", "19×42=798"),
    @("58×47=2726", "85×63=5355"),
    @("95×71=6745", "68×35=2380"),
    @("57×98=5586", "23×34=782"),
    @("13×44=572", "48×87=4176"),
    @("16×88=1408", "23×80=1840"),
    @("33×14=462", "25×45=1125"),
    @("87×41=3567", "64×19=1216"),
    @("25×91=2275", "29×38=1102"),
    @("46×46=2116", "58×37=2146"),
    @("98×32=3136", "62×31=1922"),
    @("27×69=1863", "66×64=4224"),
    @("32×40=1280", "14×72=1008")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
